$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (FilesTab) query text, corrected: the "File Type" and "Breed" output
# columns (and their coalesce() lines) were removed from the RETURN clause.
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Staffordshire Bull Terrier'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# The shorter query text now wraps onto one fewer line, so row 4 shrinks to
# match row 3's height.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moves from D3 to B4.
$ws.Range("B4").Select() | Out-Null
